$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.451.92'
$ws.Range("E2").Value = '  +2.07%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.610.21'
$ws.Range("E3").Value = '  +0.83%  '

# Row 4
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.62'
$ws.Range("E5").Value = '  +7.56%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '571.75'
$ws.Range("E6").Value = '  -0.58%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("E7").Value = '  -0.46%  '

# Row 8
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.687'
$ws.Range("E9").Value = '  +1.38%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '64.19'
$ws.Range("E10").Value = '  +14.54%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.149'
$ws.Range("E11").Value = '  -0.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000284'
$ws.Range("E12").Value = '  +4.13%  '

# Row 13
$ws.Range("E13").Value = '  +5.31%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.175.18'
$ws.Range("E14").Value = '  +0.45%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.606.46'
$ws.Range("E15").Value = '  +0.46%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.32'
$ws.Range("E16").Value = '  +5.11%  '

# Row 17
$ws.Range("E17").Value = '  +0.52%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.243.14'
$ws.Range("E18").Value = '  +1.75%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.30'
$ws.Range("E19").Value = '  +1.08%  '

# Row 20
$ws.Range("E20").Value = '  +0.76%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '407.15'
$ws.Range("E21").Value = '  +1.57%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.19'
$ws.Range("E22").Value = '  -0.43%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.39'
$ws.Range("E23").Value = '  +8.85%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.07'
$ws.Range("E24").Value = '  -0.83%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.91'
$ws.Range("E25").Value = '  -0.58%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.57'
$ws.Range("E26").Value = '  +0.84%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.87'
$ws.Range("E27").Value = '  +6.69%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.28'
$ws.Range("E28").Value = '  +3.53%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '31.73'
$ws.Range("E29").Value = '  +1.75%  '

# Row 30
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.62'
$ws.Range("E30").Value = '  -0.23%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '706.96'
$ws.Range("E31").Value = '  +10.81%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.22'
$ws.Range("E32").Value = '  +0.53%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.114'
$ws.Range("E33").Value = '  -0.24%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.77'
$ws.Range("E34").Value = '  -0.35%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '42.15'
$ws.Range("E35").Value = '  -0.36%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.416'
$ws.Range("E36").Value = '  +4.02%  '

# Row 37
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.26'
$ws.Range("E38").Value = '  +9.41%  '

# Row 39
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.21'
$ws.Range("E39").Value = '  +22.54%  '

# Row 40
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0760'
$ws.Range("E40").Value = '  -1.90%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.196.06'
$ws.Range("E41").Value = '  +0.60%  '

# Row 42
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'

# Row 44
$ws.Range("E44").Value = '  -1.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0417'
$ws.Range("E45").Value = '  +0.43%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.77'
$ws.Range("E46").Value = '  +9.35%  '

# Row 47
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.84'
$ws.Range("E47").Value = '  +3.00%  '

# Row 48
$ws.Range("E48").Value = '  +0.86%  '

# Row 49
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.10'
$ws.Range("E49").Value = '  +0.74%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '139.31'
$ws.Range("E50").Value = '  -1.63%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.72'
$ws.Range("E51").Value = '  -0.58%  '
